$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-save A64's timestamp (same instant, re-serialized with full double precision)
$ws.Range("A64").Value = 44377.76845540047

# Append the newly retrieved row of data (2021-07-01 run)
$ws.Range("A65").Value = 44378.76918391869
$ws.Range("B65").Value = 78547
$ws.Range("C65").Value = 66218
$ws.Range("D65").Value = 3643
$ws.Range("E65").Value = 2127
$ws.Range("F65").Value = 1520
$ws.Range("G65").Value = 20878
$ws.Range("H65").Value = 1566
$ws.Range("I65").Value = 877
$ws.Range("J65").Value = 202
